$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 10908.9
$ws.Range("I125").Value = 682.75
$ws.Range("K125").Value = 6144.75
$ws.Range("M125").Value = -3684.75

$ws.Range("H127").Value = 592.1177
$ws.Range("I127").Value = 384.18182
$ws.Range("J127").Value = 973.3333
$ws.Range("K127").Value = 1152.54546
$ws.Range("L127").Value = 2919.9999
$ws.Range("M127").Value = 3807.45454
$ws.Range("N127").Value = -12839.9999

$ws.Range("H129").Value = 851.7288
$ws.Range("I129").Value = 315.72726
$ws.Range("J129").Value = 974.5625
$ws.Range("K129").Value = 947.18178
$ws.Range("L129").Value = 2923.6875
$ws.Range("M129").Value = 4052.81822
$ws.Range("N129").Value = -12923.6875


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24526.578
$ws.Range("I32").Value = 26689.805
$ws.Range("J32").Value = 2353.5
$ws.Range("K32").Value = 26689.805
$ws.Range("L32").Value = 2353.5
$ws.Range("M32").Value = -26402.805
$ws.Range("N32").Value = -2927.5

$ws.Range("H45").Value = 2166.6
$ws.Range("I45").Value = 2199.9285
$ws.Range("K45").Value = 2199.9285
$ws.Range("M45").Value = -1822.9285

$ws.Range("H61").Value = 6023.9565
$ws.Range("I61").Value = 4622.143
$ws.Range("J61").Value = 10484.272
$ws.Range("K61").Value = 4622.143
$ws.Range("L61").Value = 10484.272
$ws.Range("M61").Value = -4410.143
$ws.Range("N61").Value = -10908.272

$ws.Range("H74").Value = 1686.1777
$ws.Range("I74").Value = 1365.1143
$ws.Range("J74").Value = 2809.9
$ws.Range("K74").Value = 1365.1143
$ws.Range("L74").Value = 2809.9
$ws.Range("M74").Value = -491.1143
$ws.Range("N74").Value = -4557.9

$ws.Range("H77").Value = 1686.1777
$ws.Range("I77").Value = 1365.1143
$ws.Range("J77").Value = 2809.9
$ws.Range("K77").Value = 6825.5715
$ws.Range("L77").Value = 14049.5
$ws.Range("M77").Value = -2457.5715
$ws.Range("N77").Value = -22785.5

$ws.Range("H95").Value = 49900
$ws.Range("J95").Value = 49900
$ws.Range("L95").Value = 49900
$ws.Range("N95").Value = -55392

$ws.Range("H122").Value = 1674.5
$ws.Range("I122").Value = 1566
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4698
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2248
$ws.Range("N122").Value = -10900

$ws.Range("H136").Value = 6023.9565
$ws.Range("I136").Value = 4622.143
$ws.Range("J136").Value = 10484.272
$ws.Range("K136").Value = 13866.429
$ws.Range("L136").Value = 31452.816
$ws.Range("M136").Value = -11316.429
$ws.Range("N136").Value = -36552.81600000001


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6213.357
$ws.Range("I105").Value = 6577.4443
$ws.Range("J105").Value = 5558
$ws.Range("K105").Value = 6577.4443
$ws.Range("L105").Value = 5558
$ws.Range("M105").Value = -4830.4443
$ws.Range("N105").Value = -9052


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16543.611
$ws.Range("I31").Value = 20882.5
$ws.Range("J31").Value = 11120
$ws.Range("K31").Value = 20882.5
$ws.Range("L31").Value = 11120
$ws.Range("M31").Value = -20587.5
$ws.Range("N31").Value = -11710

$ws.Range("H34").Value = 16543.611
$ws.Range("I34").Value = 20882.5
$ws.Range("J34").Value = 11120
$ws.Range("K34").Value = 20882.5
$ws.Range("L34").Value = 11120
$ws.Range("M34").Value = -20680.5
$ws.Range("N34").Value = -11524

$ws.Range("H122").Value = 18939.25
$ws.Range("I122").Value = 11000
$ws.Range("J122").Value = 23702.8
$ws.Range("K122").Value = 33000
$ws.Range("L122").Value = 71108.39999999999
$ws.Range("M122").Value = -30550
$ws.Range("N122").Value = -76008.39999999999


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 540.88
$ws.Range("I15").Value = 200
$ws.Range("J15").Value = 768.13336
$ws.Range("K15").Value = 600
$ws.Range("L15").Value = 2304.40008
$ws.Range("M15").Value = -460
$ws.Range("N15").Value = -2584.40008

$ws.Range("H55").Value = 3500.8
$ws.Range("I55").Value = 2004
$ws.Range("J55").Value = 3875
$ws.Range("K55").Value = 6012
$ws.Range("L55").Value = 11625
$ws.Range("M55").Value = -5835
$ws.Range("N55").Value = -11979

$ws.Range("H131").Value = 44065.953
$ws.Range("I131").Value = 927.2778
$ws.Range("J131").Value = 238190
$ws.Range("K131").Value = 2781.8334
$ws.Range("L131").Value = 714570
$ws.Range("M131").Value = 2258.1666
$ws.Range("N131").Value = -724650

$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 3250
$ws.Range("J132").Value = 1375
$ws.Range("K132").Value = 29250
$ws.Range("L132").Value = 12375
$ws.Range("M132").Value = -26720
$ws.Range("N132").Value = -17435

$ws.Range("H133").Value = 3470.5
$ws.Range("I133").Value = 2931.5
$ws.Range("J133").Value = 4144.25
$ws.Range("K133").Value = 8794.5
$ws.Range("L133").Value = 12432.75
$ws.Range("M133").Value = -3734.5
$ws.Range("N133").Value = -22552.75


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2575.375
$ws.Range("I122").Value = 2067.2144
$ws.Range("J122").Value = 3286.8
$ws.Range("K122").Value = 6201.6432
$ws.Range("L122").Value = 9860.400000000001
$ws.Range("M122").Value = -3751.6432
$ws.Range("N122").Value = -14760.4

$ws.Range("H139").Value = 31231.25
$ws.Range("J139").Value = 31231.25
$ws.Range("L139").Value = 31231.25
$ws.Range("N139").Value = -41511.25


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3175.6
$ws.Range("I40").Value = 2695.7778
$ws.Range("K40").Value = 2695.7778
$ws.Range("M40").Value = -2559.7778

$ws.Range("H122").Value = 6100.174
$ws.Range("I122").Value = 6005.3423
$ws.Range("K122").Value = 18016.0269
$ws.Range("M122").Value = -15566.0269


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 14521.25
$ws.Range("J58").Value = 19000
$ws.Range("L58").Value = 19000
$ws.Range("N58").Value = -19616

$ws.Range("H74").Value = 20572.111
$ws.Range("J74").Value = 23592.715
$ws.Range("L74").Value = 23592.715
$ws.Range("N74").Value = -25464.715

$ws.Range("H77").Value = 20572.111
$ws.Range("J77").Value = 23592.715
$ws.Range("L77").Value = 70778.145
$ws.Range("N77").Value = -80138.145

$ws.Range("H122").Value = 12532.777
$ws.Range("I122").Value = 1370
$ws.Range("K122").Value = 4110
$ws.Range("M122").Value = -1660

